$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "24.661.98"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.687.99"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.74%  "

# Row 5 - BNB
Set-TextValue "D5" "315.73"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.92%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3940"
$ws.Range("E7").Value = "  -0.20%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.4046"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9 - Polygon
Set-TextValue "D9" "1.486"
$ws.Range("E9").Value = "  -2.08%  "

# Row 10 - BinanceUSD
Set-TextValue "D10" "1.003"
$ws.Range("E10").Value = "  +0.83%  "

# Row 11 - OKB
Set-TextValue "D11" "52.80"
$ws.Range("E11").Value = "  -1.06%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.08800"
$ws.Range("E12").Value = "  +0.38%  "

# Row 13 - Polkadot
Set-TextValue "D13" "7.211"
$ws.Range("E13").Value = "  -1.21%  "

# Row 14 - Solana
Set-TextValue "D14" "23.44"
$ws.Range("E14").Value = "  +0.69%  "

# Row 15 - Chainlink
Set-TextValue "D15" "8.039"
$ws.Range("E15").Value = "  +7.71%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.00001312"
$ws.Range("E16").Value = "  -0.75%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.692.50"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18 - Litecoin
Set-TextValue "D18" "99.65"
$ws.Range("E18").Value = "  -0.90%  "

# Row 19 - TRON
Set-TextValue "D19" "0.07002"
$ws.Range("E19").Value = "  -0.20%  "

# Row 20 - Avalanche
Set-TextValue "D20" "19.43"
$ws.Range("E20").Value = "  -0.18%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.993"
$ws.Range("E21").Value = "  +3.76%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +1.02%  "

# Row 23 - Cosmos
Set-TextValue "D23" "14.27"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "24.655.83"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25 - LidoDAOToken
Set-TextValue "D25" "3.281"
$ws.Range("E25").Value = "  +10.79%  "

# Row 26 - Toncoin
Set-TextValue "D26" "2.364"
$ws.Range("E26").Value = "  +2.63%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "22.65"
$ws.Range("E27").Value = "  +1.21%  "

# Row 28 - Monero
Set-TextValue "D28" "162.56"
$ws.Range("E28").Value = "  +2.59%  "

# Row 29 - was BitcoinCash, now HuobiToken
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D29" "5.177"
$ws.Range("E29").Value = "  +1.14%  "

# Row 30 - was HuobiToken, now BitcoinCash
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D30" "135.21"
$ws.Range("E30").Value = "  +1.53%  "

# Row 31 - Filecoin
Set-TextValue "D31" "7.610"
$ws.Range("E31").Value = "  +2.49%  "

# Row 32 - WrappedliquidstakedEther2.0
$ws.Range("D32").Value = "1.882.11"

# Row 33 - Hedera
Set-TextValue "D33" "0.08557"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "1.057"
$ws.Range("E34").Value = "  -3.35%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue "D35" "7.176"
$ws.Range("E35").Value = "  -3.32%  "

# Row 36 - FraxShare
Set-TextValue "D36" "11.14"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37 - Algorand
Set-TextValue "D37" "0.2728"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38 - WEMIXTOKEN
Set-TextValue "D38" "1.895"
$ws.Range("E38").Value = "  -1.85%  "

# Row 39 - Aptos
$ws.Range("E39").Value = "  -3.19%  "

# Row 41 - VeChain
Set-TextValue "D41" "0.02709"
$ws.Range("E41").Value = "  -2.02%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "1.459"
$ws.Range("E42").Value = "  -1.04%  "

# Row 43 - TheSandbox
Set-TextValue "D43" "0.7582"
$ws.Range("E43").Value = "  -0.93%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "16.02"
$ws.Range("E44").Value = "  +4.32%  "

# Row 45 - NEARProtocol
Set-TextValue "D45" "2.583"
$ws.Range("E45").Value = "  +5.26%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.7119"
$ws.Range("E46").Value = "  -1.35%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "4.213"
$ws.Range("E47").Value = "  +1.47%  "

# Row 48 - Frax
$ws.Range("E48").Value = "  +1.04%  "

# Row 49 - was Quant, now Flow
$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
Set-TextValue "D49" "1.314"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50 - was Flow, now Quant
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "139.22"
$ws.Range("E50").Value = "  -1.64%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.07969"
$ws.Range("E51").Value = "  -0.73%  "
